$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    4274.010453041199,
    4109.967060358828,
    4105.498121553924,
    4011.262884979989,
    4011.262884979989,
    4011.262884979989,
    4011.262884979989,
    4011.262884979989,
    4011.262884979989,
    4011.262884979989,
    3911.81228068915
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 3).Value = $v
    $row++
}

$wb.Save()
